# Generate Report for Handoff
# - Status moves from "In Translation" to "Ready for handoff"
#   (Overview!E2, Overview!F2, zh-cn!C2, de-de!C2 all share this text)
# - The "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" timestamps
#   advance a bit over a minute later, reflecting the new handoff.
# - Column widths on the affected "Status"/language columns grow to fit the
#   new, longer "Ready for handoff" text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text: "In Translation" -> "Ready for handoff" ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# --- Timestamps ---
$wsOverview.Range("G2").Value = "2016-11-14 07:01:45"
$wsDeDe.Range("H2").Value = "2016-11-14 07:01:45"
$wsZhCn.Range("H2").Value = "2016-11-14 07:01:32"

# --- Column widths: widen the Status/language columns to fit the new text ---
$wsOverview.Columns.Item(5).ColumnWidth = 16.38265482584637
$wsOverview.Columns.Item(6).ColumnWidth = 16.38265482584637
$wsZhCn.Columns.Item(3).ColumnWidth = 16.38265482584637
$wsDeDe.Columns.Item(3).ColumnWidth = 16.38265482584637
